$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new "Attribute" column (O) by copying column N's formatting
# (borders/fills) for the used rows into the new column, then set the two
# new cell values: the header in O1 and the color value in O10.
$ws.Range("N1:N11").Copy()
$ws.Range("O1:O11").PasteSpecial(-4122)

$ws.Range("O1").Value = "Attribute"
$ws.Range("O10").Value = "#337ab7"

# Column O should be the same width as column N.
$ws.Columns("O").ColumnWidth = $ws.Columns("N").ColumnWidth

# Match the scroll position / selection the user ended up with.
$ws.Range("P10").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 8
